$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (rows 2-51: Coin, Link, Price, Volume(1h)).
# Only cells whose value actually changes are listed ($null = leave as-is).
# Column A (running index) is not touched by this update.
$rows = @(
    @{ Row = 2; B = $null; C = $null; D = "71.390.56"; E = "  +0.75%  " },
    @{ Row = 3; B = $null; C = $null; D = "3.839.48"; E = "  +1.08%  " },
    @{ Row = 4; B = $null; C = $null; D = $null; E = "  -0.01%  " },
    @{ Row = 5; B = $null; C = $null; D = "713.05"; E = "  +1.72%  " },
    @{ Row = 6; B = $null; C = $null; D = "173.39"; E = "  +0.19%  " },
    @{ Row = 7; B = $null; C = $null; D = "3.839.00"; E = $null },
    @{ Row = 8; B = $null; C = $null; D = $null; E = "  +0.03%  " },
    @{ Row = 9; B = $null; C = $null; D = "0.528"; E = "  +0.01%  " },
    @{ Row = 10; B = $null; C = $null; D = $null; E = "  +0.41%  " },
    @{ Row = 11; B = $null; C = $null; D = "7.33"; E = "  +0.98%  " },
    @{ Row = 12; B = $null; C = $null; D = $null; E = "  +0.29%  " },
    @{ Row = 13; B = $null; C = $null; D = "0.0000257"; E = "  +0.01%  " },
    @{ Row = 14; B = $null; C = $null; D = "36.83"; E = "  +2.12%  " },
    @{ Row = 15; B = $null; C = $null; D = "4.489.78"; E = "  +1.25%  " },
    @{ Row = 16; B = $null; C = $null; D = "3.924.98"; E = "  +3.39%  " },
    @{ Row = 17; B = $null; C = $null; D = "71.359.36"; E = "  +0.73%  " },
    @{ Row = 18; B = $null; C = $null; D = "7.27"; E = "  +1.07%  " },
    @{ Row = 19; B = $null; C = $null; D = $null; E = "  +0.39%  " },
    @{ Row = 20; B = $null; C = $null; D = "17.49"; E = "  -0.94%  " },
    @{ Row = 21; B = $null; C = $null; D = "497.59"; E = "  +3.65%  " },
    @{ Row = 22; B = $null; C = $null; D = $null; E = "  -1.99%  " },
    @{ Row = 23; B = $null; C = $null; D = "0.737"; E = "  +3.43%  " },
    @{ Row = 24; B = $null; C = $null; D = "85.46"; E = "  +1.72%  " },
    @{ Row = 25; B = $null; C = $null; D = $null; E = "  +2.00%  " },
    @{ Row = 26; B = $null; C = $null; D = "10.71"; E = "  +1.24%  " },
    @{ Row = 27; B = $null; C = $null; D = "12.20"; E = "  -0.61%  " },
    @{ Row = 28; B = "WrappedeETH"; C = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; D = "3.997.27"; E = "  +1.22%  " },
    @{ Row = 29; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "2.11"; E = "  -2.48%  " },
    @{ Row = 30; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "3.16"; E = "  +0.94%  " },
    @{ Row = 31; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.00"; E = "  -0.03%  " },
    @{ Row = 32; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "7.51"; E = "  -1.11%  " },
    @{ Row = 33; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "2.24"; E = "  -1.86%  " },
    @{ Row = 34; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "29.45"; E = "  +0.02%  " },
    @{ Row = 35; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "0.181"; E = "  -3.86%  " },
    @{ Row = 36; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "9.26"; E = "  +0.09%  " },
    @{ Row = 37; B = "RenzoRestakedETH"; C = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"; D = "3.807.02"; E = "  +1.61%  " },
    @{ Row = 38; B = "Binance-PegBSC-USD"; C = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D = "1.00"; E = "  +0.74%  " },
    @{ Row = 39; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.103"; E = "  +0.36%  " },
    @{ Row = 40; B = $null; C = $null; D = "6.04"; E = "  +0.43%  " },
    @{ Row = 41; B = "Mantle"; C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D = "1.04"; E = "  +4.99%  " },
    @{ Row = 42; B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D = "3.37"; E = "  -1.66%  " },
    @{ Row = 43; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "2.28"; E = "  +1.65%  " },
    @{ Row = 44; B = "USDe"; C = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; D = $null; E = "  +0.00%  " },
    @{ Row = 45; B = "FirstDigitalUSD"; C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D = "1.00"; E = "  +0.11%  " },
    @{ Row = 46; B = "FLOKI"; C = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"; D = "0.000320"; E = "  +0.49%  " },
    @{ Row = 47; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "164.00"; E = "  -0.45%  " },
    @{ Row = 48; B = "Bittensor"; C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D = "429.92"; E = "  +3.86%  " },
    @{ Row = 49; B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "49.00"; E = "  +0.36%  " },
    @{ Row = 50; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "8.76"; E = "  +1.38%  " },
    @{ Row = 51; B = "ONDO"; C = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"; D = "1.39"; E = "  -0.49%  " }
)

foreach ($r in $rows) {
    $row = $r.Row
    if ($null -ne $r.B) { $ws.Cells.Item($row, 2).Value = $r.B }
    if ($null -ne $r.C) { $ws.Cells.Item($row, 3).Value = $r.C }
    if ($null -ne $r.D) {
        # Price column: force text so values like "1.00" / "0.000320" / "71.390.56"
        # are stored verbatim and not reinterpreted as numbers.
        $ws.Cells.Item($row, 4).Value = "'" + $r.D
    }
    if ($null -ne $r.E) { $ws.Cells.Item($row, 5).Value = $r.E }
}
